# Generate Report for Handoff
# - Mark rows 7,8,9,10,11,14 (Priority column, E) as "ht" (handoff type)
#   on both the "zh-cn" and "de-de" localization sheets.
# - Refresh the "Latest Handoff Datetime" (zh-cn col H / de-de col H) and
#   the "Latest HO Xliff Generate Date" (Overview col G) timestamps for
#   those same report rows to reflect the new report-generation run.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 14)

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-08-23 20:22:12"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-08-23 20:22:17"
}

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-23 20:22:17"
}
